# Add a "Save" column (column H) to the s_vals sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold font + border) from the neighboring
# "sum" header cell (G1) onto the new H1 header cell before writing the
# text, so the Copy operation's value doesn't clobber what we set below.
$ws.Range("G1").Copy($ws.Range("H1"))

# Header text for the new column.
$ws.Range("H1").Value = "Save"

# Data values for the new column (Save = 0 for the first row, 1 for the
# second row), matching the rest of the unstyled numeric data cells.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
